$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Expense" to "Income"
$ws.Name = "Income"

# Header row
$ws.Range("A1").Value = "Source"
$ws.Range("B1").Value = "Amount"
$ws.Range("C1").Value = "Date"

# Data rows (Source, Amount, Date)
$ws.Range("A2").Value = "salary"
$ws.Range("B2").Value = 5000
$ws.Range("C2").Value = 45879.22928240741

$ws.Range("A3").Value = "salary"
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 45689.22928240741

$ws.Range("A4").Value = "salary"
$ws.Range("B4").Value = 500
$ws.Range("C4").Value = 45689.22928240741

$ws.Range("A5").Value = "salary"
$ws.Range("B5").Value = 5000
$ws.Range("C5").Value = 45689.22928240741

# Apply the built-in short-date number format (numFmtId 14) to C2, then
# propagate the exact same style to C3:C5 via copy/paste-format so the
# workbook ends up with a single shared date style instead of one per cell.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C3:C5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
